{"js": "// Add a new paragraph containing \"678910\" at the end of the document\n// body, right after the existing \"12345\" paragraph.\nconst body = context.document.body;\nbody.insertParagraph(\"678910\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Add a new paragraph containing \"678910\" after the existing last\n# paragraph (\"12345\") in the document body.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Collapse(0)  # wdCollapseEnd\n$rng.InsertAfter(\"`r678910\")\n"}
